$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 374.68182
$ws.Range("K92").Value = 334.6316
$ws.Range("I92").Value = 334.6316
$ws.Range("M92").Value = 913.3684000000001
$ws.Range("I98").Value = 1241.2
$ws.Range("H98").Value = 1241.2
$ws.Range("K98").Value = 1241.2
$ws.Range("M98").Value = 256.8
$ws.Range("I100").Value = 1333.4445
$ws.Range("M100").Value = -792.4445000000001
$ws.Range("K100").Value = 1333.4445
$ws.Range("J100").Value = 1800
$ws.Range("L100").Value = 1800
$ws.Range("N100").Value = -2882
$ws.Range("J112").Value = 1969.2667
$ws.Range("N112").Value = -8123.800099999999
$ws.Range("H112").Value = 1896.1875
$ws.Range("L112").Value = 5907.800099999999
$ws.Range("K122").Value = 3723.6
$ws.Range("M122").Value = -1273.6
$ws.Range("I122").Value = 1241.2
$ws.Range("H122").Value = 1241.2
$ws.Range("M129").Value = 3245
$ws.Range("I129").Value = 585
$ws.Range("J129").Value = 2637.75
$ws.Range("L129").Value = 7913.25
$ws.Range("N129").Value = -17913.25
$ws.Range("H129").Value = 1881.4736
$ws.Range("K129").Value = 1755
$ws.Range("K132").Value = 6781.38
$ws.Range("M132").Value = -4251.38
$ws.Range("H132").Value = 2522.8154
$ws.Range("I132").Value = 2260.46
$ws.Range("I137").Value = 1453.3077
$ws.Range("H137").Value = 29712.834
$ws.Range("N137").Value = -314662.8
$ws.Range("K137").Value = 4359.9231
$ws.Range("J137").Value = 103187.6
$ws.Range("L137").Value = 309562.8
$ws.Range("M137").Value = -1809.9231
$ws.Range("L138").Value = 16945.2348
$ws.Range("H138").Value = 3016.95
$ws.Range("I138").Value = 1071.9565
$ws.Range("J138").Value = 5648.4116
$ws.Range("N138").Value = -27225.2348
$ws.Range("K138").Value = 3215.8695
$ws.Range("M138").Value = 1924.1305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N62").Value = -31248
$ws.Range("L62").Value = 30000
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("J65").Value = 30000
$ws.Range("N65").Value = -96240
$ws.Range("H65").Value = 30000
$ws.Range("K74").Value = 1740.4572
$ws.Range("J74").Value = 1779.4166
$ws.Range("H74").Value = 1750.4043
$ws.Range("N74").Value = -3527.4166
$ws.Range("I74").Value = 1740.4572
$ws.Range("M74").Value = -866.4572000000001
$ws.Range("L74").Value = 1779.4166
$ws.Range("J77").Value = 1779.4166
$ws.Range("M77").Value = -4334.286
$ws.Range("L77").Value = 8897.083000000001
$ws.Range("N77").Value = -17633.083
$ws.Range("H77").Value = 1750.4043
$ws.Range("K77").Value = 8702.286
$ws.Range("I77").Value = 1740.4572
$ws.Range("L102").Value = 2118
$ws.Range("J102").Value = 2118
$ws.Range("K102").Value = 1053.5625
$ws.Range("N102").Value = -5362
$ws.Range("M102").Value = 568.4375
$ws.Range("I102").Value = 1053.5625
$ws.Range("H102").Value = 1266.45
$ws.Range("N122").Value = -9250
$ws.Range("K122").Value = 2144.4999
$ws.Range("M122").Value = 305.5001000000002
$ws.Range("I122").Value = 714.8333
$ws.Range("J122").Value = 1450
$ws.Range("L122").Value = 4350
$ws.Range("H122").Value = 771.38464
$ws.Range("K132").Value = 3641.25
$ws.Range("M132").Value = -1111.25
$ws.Range("J132").Value = 1806.1428
$ws.Range("L132").Value = 5418.428400000001
$ws.Range("H132").Value = 1417.6885
$ws.Range("I132").Value = 1213.75
$ws.Range("N132").Value = -10478.4284
$ws.Range("N139").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K99").Value = 20693.8
$ws.Range("H99").Value = 8853.467000000001
$ws.Range("M99").Value = -19195.8
$ws.Range("L99").Value = 2933.3
$ws.Range("J99").Value = 2933.3
$ws.Range("N99").Value = -5929.3
$ws.Range("I99").Value = 20693.8
$ws.Range("M134").Value = -666.1032
$ws.Range("I134").Value = 1067.0344
$ws.Range("H134").Value = 1292.2253
$ws.Range("K134").Value = 3201.1032
$ws.Range("H137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("L138").Value = 63818.75
$ws.Range("H138").Value = 63818.75
$ws.Range("J138").Value = 63818.75
$ws.Range("N138").Value = -74098.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K22").Value = 324.8
$ws.Range("N22").Value = -1083.33334
$ws.Range("M22").Value = 25.19999999999999
$ws.Range("I22").Value = 324.8
$ws.Range("L22").Value = 383.33334
$ws.Range("J22").Value = 383.33334
$ws.Range("H22").Value = 338.30768
$ws.Range("N31").Value = -4790.6313
$ws.Range("L31").Value = 4200.6313
$ws.Range("I31").Value = 1514.6666
$ws.Range("K31").Value = 1514.6666
$ws.Range("M31").Value = -1219.6666
$ws.Range("H31").Value = 2496.077
$ws.Range("J31").Value = 4200.6313
$ws.Range("M34").Value = -1312.6666
$ws.Range("L34").Value = 4200.6313
$ws.Range("J34").Value = 4200.6313
$ws.Range("I34").Value = 1514.6666
$ws.Range("H34").Value = 2496.077
$ws.Range("N34").Value = -4604.6313
$ws.Range("K34").Value = 1514.6666
$ws.Range("K132").Value = 2864.4375
$ws.Range("M132").Value = -334.4375
$ws.Range("J132").Value = 2234.5293
$ws.Range("L132").Value = 6703.5879
$ws.Range("H132").Value = 1289.5077
$ws.Range("I132").Value = 954.8125
$ws.Range("N132").Value = -11763.5879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I68").Value = 416.33334
$ws.Range("J68").Value = 600
$ws.Range("H68").Value = 508.16666
$ws.Range("K68").Value = 1249.00002
$ws.Range("L68").Value = 1800
$ws.Range("N68").Value = -3422
$ws.Range("M68").Value = -438.0000199999999
$ws.Range("M71").Value = 308.9999399999997
$ws.Range("N71").Value = -13512
$ws.Range("I71").Value = 416.33334
$ws.Range("L71").Value = 5400
$ws.Range("J71").Value = 600
$ws.Range("K71").Value = 3747.00006
$ws.Range("H71").Value = 508.16666
$ws.Range("K107").Value = 975
$ws.Range("J107").Value = 552.5454999999999
$ws.Range("I107").Value = 325
$ws.Range("M107").Value = 945
$ws.Range("L107").Value = 1657.6365
$ws.Range("N107").Value = -5497.6365
$ws.Range("H107").Value = 491.86667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K102").Value = 1127.3572
$ws.Range("M102").Value = 494.6428000000001
$ws.Range("I102").Value = 1127.3572
$ws.Range("H102").Value = 1127.3572
$ws.Range("K132").Value = 12362.118
$ws.Range("M132").Value = -9832.118
$ws.Range("J132").Value = 3890.5
$ws.Range("L132").Value = 11671.5
$ws.Range("H132").Value = 4068.3865
$ws.Range("I132").Value = 4120.706
$ws.Range("N132").Value = -16731.5
$ws.Range("L135").Value = 43943.5
$ws.Range("H135").Value = 43943.5
$ws.Range("N135").Value = -54083.5
$ws.Range("J135").Value = 43943.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K132").Value = 4917.8043
$ws.Range("M132").Value = -2387.8043
$ws.Range("J132").Value = 706.6667
$ws.Range("L132").Value = 2120.0001
$ws.Range("H132").Value = 1611.29
$ws.Range("I132").Value = 1639.2681
$ws.Range("N132").Value = -7180.0001
$ws.Range("L140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("H140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M96").ClearContents()
$ws.Range("I96").Value = 0
$ws.Range("L96").Value = 2130
$ws.Range("J96").Value = 2130
$ws.Range("N96").Value = -4876
$ws.Range("H96").Value = 2130
$ws.Range("K96").Value = 0
$ws.Range("K132").Value = 3856.35
$ws.Range("M132").Value = -1326.35
$ws.Range("J132").Value = 2927.2632
$ws.Range("L132").Value = 8781.7896
$ws.Range("H132").Value = 2085.3076
$ws.Range("I132").Value = 1285.45
$ws.Range("N132").Value = -13841.7896
$ws.Range("H136").Value = 3111.8823
$ws.Range("K136").Value = 9051.293099999999
$ws.Range("N136").Value = -15601.5
$ws.Range("I136").Value = 3017.0977
$ws.Range("M136").Value = -6501.293099999999
$ws.Range("J136").Value = 3500.5
$ws.Range("L136").Value = 10501.5
